# Edit resource() / add resource() / delete resource() test data:
#  - Row 2, column C (title): "Новий ресурс для DAO у верхньому меню" -> "Новий ресурс вверху "
#  - Row 2, column F (addTextToTitle): "Ресурс відредактований" -> "відредагований"
#  - Row 3, column C (title): 'Новий ресурс для DAO в розділі "Ресурси"' -> "Новий ресурс в меню "
#  - Row 3, column F (addTextToTitle): "Ресурс відредактований" -> "відредагований"
#  - Active selection moves from G1 to F3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "addTextToTitle" column first so the new shared string
# "відредагований" is registered before the other two new strings.
$ws.Range("F2").Value = "відредагований"
$ws.Range("F3").Value = "відредагований"

$ws.Range("C2").Value = "Новий ресурс вверху "
$ws.Range("C3").Value = "Новий ресурс в меню "

$ws.Range("F3").Select()
